$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "1" to "ქ. თბილისი"
$ws.Name = "ქ. თბილისი"

# Remove the "(census results)" note from A2, leaving the row blank
# (this also drops the now-unused shared string automatically)
$ws.Range("A2").Clear()

# Remove the blank spacer row (old row 3), shifting the rows below up by one
$ws.Rows("3:3").Delete()

# Drop the historical 1989 and 2002 columns, keeping only the 2014 figures;
# the remaining 2014 column shifts left into column B
$ws.Range("B5:C5").EntireColumn.Delete()

# Restore the original selection shown in the saved file
[void]$ws.Range("A2").Select()
